# This script refreshes cached Leve market-price and profit figures
# (columns H:N -- currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ) across all eight crafting-job sheets, matching the
# latest scheduled market-data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1615.5807
$ws.Range("I132").Value = 1543.1333
$ws.Range("J132").Value = 3789
$ws.Range("K132").Value = 4629.3999
$ws.Range("L132").Value = 11367
$ws.Range("M132").Value = -2099.3999
$ws.Range("N132").Value = -16427
$ws.Range("H133").Value = 50041.8
$ws.Range("J133").Value = 49875
$ws.Range("L133").Value = 49875
$ws.Range("N133").Value = -59995
$ws.Range("H135").Value = 3184.077
$ws.Range("I135").Value = 1841.762
$ws.Range("J135").Value = 8821.799999999999
$ws.Range("K135").Value = 16575.858
$ws.Range("L135").Value = 79396.2
$ws.Range("M135").Value = -14040.858
$ws.Range("N135").Value = -84466.2
$ws.Range("H137").Value = 967.9697
$ws.Range("I137").Value = 773.15
$ws.Range("J137").Value = 1267.6923
$ws.Range("K137").Value = 2319.45
$ws.Range("L137").Value = 3803.0769
$ws.Range("M137").Value = 230.5500000000002
$ws.Range("N137").Value = -8903.0769
$ws.Range("H138").Value = 2115.2334
$ws.Range("I138").Value = 1626.909
$ws.Range("J138").Value = 3458.125
$ws.Range("K138").Value = 4880.727000000001
$ws.Range("L138").Value = 10374.375
$ws.Range("M138").Value = 259.2729999999992
$ws.Range("N138").Value = -20654.375
$ws.Range("H139").Value = 68860
$ws.Range("J139").Value = 68860
$ws.Range("L139").Value = 68860
$ws.Range("N139").Value = -79140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3867.11
$ws.Range("I32").Value = 3022.3604
$ws.Range("J32").Value = 9056.286
$ws.Range("K32").Value = 3022.3604
$ws.Range("L32").Value = 9056.286
$ws.Range("M32").Value = -2735.3604
$ws.Range("N32").Value = -9630.286
$ws.Range("H61").Value = 3045.6572
$ws.Range("I61").Value = 3491.2
$ws.Range("J61").Value = 1931.8
$ws.Range("K61").Value = 3491.2
$ws.Range("L61").Value = 1931.8
$ws.Range("M61").Value = -3279.2
$ws.Range("N61").Value = -2355.8
$ws.Range("H74").Value = 2534.2727
$ws.Range("I74").Value = 2351
$ws.Range("J74").Value = 2927
$ws.Range("K74").Value = 2351
$ws.Range("L74").Value = 2927
$ws.Range("M74").Value = -1477
$ws.Range("N74").Value = -4675
$ws.Range("H77").Value = 2534.2727
$ws.Range("I77").Value = 2351
$ws.Range("J77").Value = 2927
$ws.Range("K77").Value = 11755
$ws.Range("L77").Value = 14635
$ws.Range("M77").Value = -7387
$ws.Range("N77").Value = -23371
$ws.Range("H132").Value = 1863.0927
$ws.Range("I132").Value = 1610.2307
$ws.Range("K132").Value = 4830.6921
$ws.Range("M132").Value = -2300.6921
$ws.Range("H133").Value = 82104.39999999999
$ws.Range("J133").Value = 82104.39999999999
$ws.Range("L133").Value = 82104.39999999999
$ws.Range("N133").Value = -87164.39999999999
$ws.Range("H136").Value = 3045.6572
$ws.Range("I136").Value = 3491.2
$ws.Range("J136").Value = 1931.8
$ws.Range("K136").Value = 10473.6
$ws.Range("L136").Value = 5795.4
$ws.Range("M136").Value = -7923.599999999999
$ws.Range("N136").Value = -10895.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 937.8125
$ws.Range("J94").Value = 1065.7693
$ws.Range("L94").Value = 1065.7693
$ws.Range("N94").Value = -1967.7693
$ws.Range("H132").Value = 59800
$ws.Range("J132").Value = 59800
$ws.Range("L132").Value = 59800
$ws.Range("N132").Value = -69920
$ws.Range("H138").Value = 58920
$ws.Range("J138").Value = 58920
$ws.Range("L138").Value = 58920
$ws.Range("N138").Value = -69200
$ws.Range("H140").Value = 63966.668
$ws.Range("J140").Value = 63966.668
$ws.Range("L140").Value = 63966.668
$ws.Range("N140").Value = -74326.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4160.3716
$ws.Range("I31").Value = 3335.6785
$ws.Range("J31").Value = 7459.143
$ws.Range("K31").Value = 3335.6785
$ws.Range("L31").Value = 7459.143
$ws.Range("M31").Value = -3040.6785
$ws.Range("N31").Value = -8049.143
$ws.Range("H34").Value = 4160.3716
$ws.Range("I34").Value = 3335.6785
$ws.Range("J34").Value = 7459.143
$ws.Range("K34").Value = 3335.6785
$ws.Range("L34").Value = 7459.143
$ws.Range("M34").Value = -3133.6785
$ws.Range("N34").Value = -7863.143
$ws.Range("H58").Value = 76924240
$ws.Range("I58").Value = 142857740
$ws.Range("J58").Value = 1813.3334
$ws.Range("K58").Value = 142857740
$ws.Range("L58").Value = 1813.3334
$ws.Range("M58").Value = -142857537
$ws.Range("N58").Value = -2219.3334
$ws.Range("H134").Value = 7408547
$ws.Range("I134").Value = 9091707
$ws.Range("J134").Value = 2640
$ws.Range("K134").Value = 27275121
$ws.Range("L134").Value = 7920
$ws.Range("M134").Value = -27272586
$ws.Range("N134").Value = -12990
$ws.Range("H136").Value = 76924240
$ws.Range("I136").Value = 142857740
$ws.Range("J136").Value = 1813.3334
$ws.Range("K136").Value = 428573220
$ws.Range("L136").Value = 5440.0002
$ws.Range("M136").Value = -428570670
$ws.Range("N136").Value = -10540.0002
$ws.Range("H138").Value = 49950
$ws.Range("J138").Value = 49950
$ws.Range("L138").Value = 49950
$ws.Range("N138").Value = -60230
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1233.3334
$ws.Range("J5").Value = 966.6667
$ws.Range("L5").Value = 2900.0001
$ws.Range("N5").Value = -3124.0001
$ws.Range("H131").Value = 14286796
$ws.Range("I131").Value = 1956.6666
$ws.Range("J131").Value = 17242280
$ws.Range("K131").Value = 5869.9998
$ws.Range("L131").Value = 51726840
$ws.Range("M131").Value = -829.9997999999996
$ws.Range("N131").Value = -51736920
$ws.Range("H132").Value = 1122.1904
$ws.Range("I132").Value = 1657.8
$ws.Range("J132").Value = 954.8125
$ws.Range("K132").Value = 14920.2
$ws.Range("L132").Value = 8593.3125
$ws.Range("M132").Value = -12390.2
$ws.Range("N132").Value = -13653.3125
$ws.Range("H135").Value = 1233.3334
$ws.Range("J135").Value = 966.6667
$ws.Range("L135").Value = 8700.0003
$ws.Range("N135").Value = -13770.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 59710
$ws.Range("J133").Value = 59710
$ws.Range("L133").Value = 59710
$ws.Range("N133").Value = -69830
$ws.Range("H140").Value = 96966.664
$ws.Range("J140").Value = 96966.664
$ws.Range("L140").Value = 96966.664
$ws.Range("N140").Value = -107326.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 90764.28999999999
$ws.Range("J133").Value = 90764.28999999999
$ws.Range("L133").Value = 90764.28999999999
$ws.Range("N133").Value = -95824.28999999999
$ws.Range("H136").Value = 3017.814
$ws.Range("I136").Value = 2123.1155
$ws.Range("J136").Value = 4386.1763
$ws.Range("K136").Value = 6369.3465
$ws.Range("L136").Value = 13158.5289
$ws.Range("M136").Value = -3819.3465
$ws.Range("N136").Value = -18258.5289
$ws.Range("H139").Value = 55920
$ws.Range("J139").Value = 55920
$ws.Range("L139").Value = 55920
$ws.Range("N139").Value = -66200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58825412
$ws.Range("I81").Value = 1809
$ws.Range("J81").Value = 125001970
$ws.Range("K81").Value = 3618
$ws.Range("L81").Value = 250003940
$ws.Range("M81").Value = -2557
$ws.Range("N81").Value = -250006062
$ws.Range("H84").Value = 58825412
$ws.Range("I84").Value = 1809
$ws.Range("J84").Value = 125001970
$ws.Range("K84").Value = 18090
$ws.Range("L84").Value = 1250019700
$ws.Range("M84").Value = -12786
$ws.Range("N84").Value = -1250030308
$ws.Range("H106").Value = 3000
$ws.Range("I106").Value = 3000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1738
$ws.Range("N106").ClearContents()
$ws.Range("H133").Value = 32207.5
$ws.Range("J133").Value = 32207.5
$ws.Range("L133").Value = 32207.5
$ws.Range("N133").Value = -42327.5
$ws.Range("H136").Value = 969.8077
$ws.Range("I136").Value = 682.5
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 2047.5
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = 502.5
$ws.Range("N136").Value = -12750
$ws.Range("H138").Value = 47011.11
$ws.Range("J138").Value = 47011.11
$ws.Range("L138").Value = 47011.11
$ws.Range("N138").Value = -57291.11
$ws.Range("H141").Value = 55866.668
$ws.Range("J141").Value = 56850
$ws.Range("L141").Value = 56850
$ws.Range("N141").Value = -67210
